# Insert a new data row at row 112 (pushing the existing rows 112-183 down
# to 113-184) and populate it with the new "Femacal de La Calera" / Pepino
# ensalada observation dated 2021-08-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(112).EntireRow.Insert()

$ws.Range("A112").Value = 3
$ws.Range("B112").Value = "Femacal de La Calera"
$ws.Range("C112").Value = "Coquimbo"
$ws.Range("D112").Value = 44438
$ws.Range("E112").Value = 5
$ws.Range("F112").Value = 100112043
$ws.Range("G112").Value = "Pepino ensalada"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 140
$ws.Range("K112").Value = 11000
$ws.Range("L112").Value = 12000
$ws.Range("M112").Value = 11500
$ws.Range("N112").Value = "$/caja 70 unidades"
$ws.Range("O112").Value = "Región de Arica y Parinacota"
$ws.Range("P112").Value = 164
$ws.Range("Q112").Value = 70
$ws.Range("R112").Value = "Hortaliza"
